# Auto-generated edit script applying cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.718.10"
$ws.Range("E2").Value = "  -1.97%  "
$ws.Range("D3").Value = "2.417.77"
$ws.Range("E3").Value = "  +5.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.565"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.516"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.19%  "
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.53%  "
$ws.Range("E13").Value = "  +1.03%  "
$ws.Range("D14").Value = "2.771.41"
$ws.Range("E14").Value = "  +4.58%  "
$ws.Range("D15").Value = "2.414.21"
$ws.Range("E15").Value = "  +4.94%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.850"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.40%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.42%  "
$ws.Range("D18").Value = "45.684.95"
$ws.Range("E18").Value = "  -2.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.30%  "
$ws.Range("E20").Value = "  +1.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "243.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.62%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "39.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.85%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("B30").Value = "LidoDAOToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.87"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +17.93%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "21.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.21%  "
$ws.Range("E32").Value = "  -1.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "148.59"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.92%  "
$ws.Range("E35").Value = "  -2.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +13.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.114"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("E38").Value = "  -1.45%  "
$ws.Range("E39").Value = "  -4.48%  "
$ws.Range("E40").Value = "  -2.18%  "
$ws.Range("E41").Value = "  -0.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.29"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.37%  "
$ws.Range("D43").Value = "1.955.75"
$ws.Range("E43").Value = "  +7.57%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -10.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.22%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.42%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "15.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +15.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.189"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.24%  "
$ws.Range("D51").Value = "2.655.39"
$ws.Range("E51").Value = "  +5.11%  "
